$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '37.225.66'
$ws.Range('E2').Value2 = '  +0.31%  '
$ws.Range('D3').Value2 = '2.070.09'
$ws.Range('E3').Value2 = '  +3.89%  '
$ws.Range('E4').Value2 = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '235.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  -3.03%  '
$ws.Range('E6').Value2 = '  +2.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '57.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  +4.91%  '
$ws.Range('E9').Value2 = '  +2.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '58.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.0759'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  +1.14%  '
$ws.Range('E12').Value2 = '  +3.05%  '
$ws.Range('D13').Value2 = '2.373.65'
$ws.Range('E13').Value2 = '  +3.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '14.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +3.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '21.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  +0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '0.774'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  +2.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '5.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  +3.68%  '
$ws.Range('D18').Value2 = '2.060.21'
$ws.Range('E18').Value2 = '  +3.38%  '
$ws.Range('D19').Value2 = '37.161.16'
$ws.Range('E19').Value2 = '  +0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '5.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  +19.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '68.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  +0.27%  '
$ws.Range('D22').Value2 = '0.0₃0809'
$ws.Range('E22').Value2 = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '223.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  -1.41%  '
$ws.Range('E24').Value2 = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '2.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +2.08%  '
$ws.Range('E26').Value2 = '  +0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '162.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  +0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '8.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +2.24%  '
$ws.Range('E29').Value2 = '  +6.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '19.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +0.98%  '
$ws.Range('E31').Value2 = '  +7.51%  '
$ws.Range('E32').Value2 = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '4.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  +1.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.0621'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +1.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '2.52'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  +6.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '4.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  +3.86%  '
$ws.Range('E37').Value2 = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '5.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  +13.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '3.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  +0.80%  '
$ws.Range('E40').Value2 = '  -0.75%  '
$ws.Range('E41').Value2 = '  -2.35%  '
$ws.Range('B42').Value2 = 'Cronos'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.0965'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  +8.83%  '
$ws.Range('B43').Value2 = 'FTXToken'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '4.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  +23.91%  '
$ws.Range('D44').Value2 = '1.472.16'
$ws.Range('E44').Value2 = '  +3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '94.63'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +7.16%  '
$ws.Range('E46').Value2 = '  +2.57%  '
$ws.Range('B47').Value2 = 'InjectiveProtocol'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '16.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  +5.07%  '
$ws.Range('B48').Value2 = 'TrustWalletToken'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '1.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  +1.00%  '
$ws.Range('E49').Value2 = '  +2.14%  '
$ws.Range('E50').Value2 = '  +7.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '2.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  +1.86%  '
